$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("algorithm")

# Row 1: H1 becomes a plain number; I1 picks up the (renamed) shared string
# "cutoff=2" that used to live in H1 ("304,205 cterms").
$ws.Range("H1").Value = 148039
$ws.Range("I1").Value = "cutoff=2"

# Row 2: updated G/H values
$ws.Range("G2").Value = 1.5
$ws.Range("H2").Value = 0.032638888888888891

# Row 3: newly added G/H values (H3 uses the elapsed-time [h]:mm:ss format,
# matching the existing style used on sheet "indexing" column F)
$ws.Range("G3").Value = 17.399999999999999
$ws.Range("H3").Value = 1.5131944444444445
$ws.Range("H3").NumberFormat = "[h]:mm:ss"

# Row 4: updated G/H values
$ws.Range("G4").Value = 5.98
$ws.Range("H4").Value = 0.74444444444444446

# Row 5: newly added G/H values (H5 uses the h:mm format already used by
# H2/H4/C2..C12 on this sheet)
$ws.Range("G5").Value = 2.72
$ws.Range("H5").Value = 0.081944444444444445
$ws.Range("H5").NumberFormat = "h:mm"
